$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set shared_links ("hiatus") for rows 4-6 and update the combined column accordingly
$ws.Range("B4").Value = "hiatus"
$ws.Range("D4").Value = "Service providers performance: hiatus"

$ws.Range("B5").Value = "hiatus"
$ws.Range("D5").Value = "Investment status: hiatus"

$ws.Range("B6").Value = "hiatus"
$ws.Range("D6").Value = "Lessons Learned: hiatus"
